$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.894.24"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "2.230.47"
$ws.Range("E3").Value = "  -4.86%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'293.54"
$ws.Range("E5").Value = "  -5.40%  "
$ws.Range("D6").Value = "'85.02"
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("D7").Value = "'0.513"
$ws.Range("E7").Value = "  -2.37%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.469"
$ws.Range("E9").Value = "  -2.52%  "
$ws.Range("D10").Value = "'0.0799"
$ws.Range("E10").Value = "  -0.81%  "
$ws.Range("D11").Value = "'30.24"
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("D12").Value = "'47.94"
$ws.Range("E12").Value = "  -8.67%  "
$ws.Range("E13").Value = "  -2.30%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'6.35"
$ws.Range("E14").Value = "  -0.76%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.572.82"
$ws.Range("E15").Value = "  -4.98%  "
$ws.Range("D16").Value = "'14.19"
$ws.Range("E16").Value = "  -4.26%  "
$ws.Range("D17").Value = "2.225.39"
$ws.Range("E17").Value = "  -5.61%  "
$ws.Range("D18").Value = "'0.723"
$ws.Range("E18").Value = "  -4.55%  "
$ws.Range("D19").Value = "39.806.25"
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("D20").Value = "0.0₃0891"
$ws.Range("E20").Value = "  -1.05%  "
$ws.Range("D21").Value = "'5.79"
$ws.Range("E21").Value = "  -4.80%  "
$ws.Range("D22").Value = "'65.46"
$ws.Range("E22").Value = "  -3.91%  "
$ws.Range("D23").Value = "'10.54"
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("E24").Value = "  -1.18%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  -4.99%  "
$ws.Range("E27").Value = "  +1.25%  "
$ws.Range("D28").Value = "'22.95"
$ws.Range("E28").Value = "  -2.74%  "
$ws.Range("E29").Value = "  +2.90%  "
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("D31").Value = "'154.86"
$ws.Range("E31").Value = "  +0.81%  "
$ws.Range("D32").Value = "'32.90"
$ws.Range("E32").Value = "  -5.80%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  -0.30%  "
$ws.Range("E34").Value = "  -4.73%  "
$ws.Range("E35").Value = "  -1.37%  "
$ws.Range("E36").Value = "  -5.22%  "
$ws.Range("D37").Value = "'16.44"
$ws.Range("E37").Value = "  +5.87%  "
$ws.Range("E38").Value = "  -1.59%  "
$ws.Range("D39").Value = "'0.0981"
$ws.Range("E39").Value = "  -0.79%  "
$ws.Range("D40").Value = "'2.67"
$ws.Range("E40").Value = "  -3.75%  "
$ws.Range("E41").Value = "  -3.09%  "
$ws.Range("E42").Value = "  -3.07%  "
$ws.Range("D43").Value = "1.950.47"
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("E44").Value = "  -3.06%  "
$ws.Range("E45").Value = "  +1.73%  "
$ws.Range("D46").Value = "'9.37"
$ws.Range("E46").Value = "  -0.72%  "
$ws.Range("D47").Value = "'16.26"
$ws.Range("E47").Value = "  -6.63%  "
$ws.Range("E48").Value = "  -2.62%  "
$ws.Range("D49").Value = "2.444.37"
$ws.Range("E49").Value = "  -4.83%  "
$ws.Range("D50").Value = "'70.86"
$ws.Range("E50").Value = "  +0.99%  "
$ws.Range("D51").Value = "'88.88"
$ws.Range("E51").Value = "  -4.38%  "
